$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (48) with the latest bitcoin buy entry dated 09/14/2025.
# Force column A to be treated as text so the date string isn't converted
# into a numeric date serial value, then restore the default/general
# number format so no extra cell style is left behind.
$ws.Cells.Item(48, 1).NumberFormat = "@"
$ws.Cells.Item(48, 1).Value = "09/14/2025"
$ws.Cells.Item(48, 1).ClearFormats()
$ws.Cells.Item(48, 2).Value = 0.0004291899999999994
$ws.Cells.Item(48, 3).Value = 116498.5204687902
$ws.Cells.Item(48, 4).Value = 50
